$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily-scrape row lands at row 893 ("2026/02/28  土  16  201"),
# pushing the existing rows 893:934 down to 894:935.
$ws.Rows(893).Insert()

# Column A holds date text (not a real Excel date), so force Text format
# before assigning, otherwise COM auto-coerces "2026/02/28" into a date
# serial. ClearFormats() afterwards drops the now-unneeded explicit
# number-format style, leaving the cell on the sheet's default style -
# matching how the rest of the column is stored.
$cellA = $ws.Cells.Item(893, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026/02/28"
$cellA.ClearFormats()

$ws.Cells.Item(893, 2).Value = "土"
$ws.Cells.Item(893, 3).Value = 16
$ws.Cells.Item(893, 4).Value = 201
